# Generate Report for Handback
# The file "dfb2261e-ee4c-474d-9114-347e53dcc0ad.md" has now been handed
# back (was previously "Ready for handoff"). Update its status on the
# Overview sheet and on each locale sheet, and stamp the new handback
# datetime for each locale.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# Overview sheet: both locale status columns for this file move to
# "Handed back: in sync with en-US".
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $statusHandedBack
$overview.Range("C3").Value = $statusHandedBack

# zh-cn detail sheet: Status column (C) + Latest Handback DateTime (H)
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $statusHandedBack

# de-de detail sheet: Status column (C) + Latest Handback DateTime (H)
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $statusHandedBack

# Stamp the new handback timestamps recorded for this handback run.
$zhcn.Range("H3").Value = "2016-03-18 08:35:29"
$dede.Range("H3").Value = "2016-03-18 08:35:34"
